$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: a new record (week of 2023-08-09) is inserted at the
# top of the "Vega Modelo de Temuco - Mango" data block (row 599), pushing
# all the existing rows from 599..625 down by one (they become 600..626).
$ws.Rows.Item(599).Insert()

$ws.Cells.Item(599, 1).Value = 10
$ws.Cells.Item(599, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(599, 3).Value = "La Araucanía"
$ws.Cells.Item(599, 4).Value = 45147
$ws.Cells.Item(599, 5).Value = 9
$ws.Cells.Item(599, 6).Value = "Fruta"
$ws.Cells.Item(599, 7).Value = 100108
$ws.Cells.Item(599, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(599, 9).Value = 100108002
$ws.Cells.Item(599, 10).Value = "Mango"
$ws.Cells.Item(599, 11).Value = "Sin especificar"
$ws.Cells.Item(599, 12).Value = "Primera"
$ws.Cells.Item(599, 13).Value = 600
$ws.Cells.Item(599, 14).Value = 8500
$ws.Cells.Item(599, 15).Value = 9000
$ws.Cells.Item(599, 16).Value = 8750
$ws.Cells.Item(599, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(599, 18).Value = "Brasil"
$ws.Cells.Item(599, 19).Value = 2188
$ws.Cells.Item(599, 20).Value = 4
